# [Fonds de solidarite] Add 2020-12-09 data
# Update nombre_aides (C) and montant_total (D) figures for several
# region / classe_effectif rows to reflect the 2020-12-09 data refresh.
# Values are entered with a leading apostrophe so Excel keeps them as text
# (matching the source data, which stores these columns as text/inline
# strings rather than numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   C = "1385"; D = "6376555.24" },
    @{ Row = 4;   C = "598";  D = "4952932.87" },
    @{ Row = 25;  C = "73";   D = "935024.25" },
    @{ Row = 56;  C = "1010"; D = "5588551.94" },
    @{ Row = 57;  C = "486";  D = "4333391.40" },
    @{ Row = 58;  C = "177";  D = "1700340.18" },
    @{ Row = 59;  C = "67";   D = "995420.00" },
    @{ Row = 60;  C = "14";   D = "301032.25" },
    @{ Row = 63;  C = "5688"; D = "23860987.63" },
    @{ Row = 64;  C = "3108"; D = "18466223.98" },
    @{ Row = 105; C = "495";  D = "2305228.10" }
)

foreach ($u in $updates) {
    $cellC = $ws.Cells.Item($u.Row, 3)
    $cellC.Value = "'" + $u.C
    $cellC.Style = "Normal"

    $cellD = $ws.Cells.Item($u.Row, 4)
    $cellD.Value = "'" + $u.D
    $cellD.Style = "Normal"
}
